$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "exercises/e03.html"
$ws.Range("E4").Value = "slides/slides.html#/sitzung-03-eine-problemstellung-entwickeln"
$ws.Range("E3").Value = "slides/slides.html#/sitzung-02-parasoziale-beziehungen-im-zeitverlauf"
$ws.Range("E5").Value = "slides/slides.html#/sitzung-04-grundlagen-der-manuellen-inhaltsanalyse"

$ws.Range("E5").Select()
